# Update the date heading and the 25 division-problem table cells.
$d = $word.ActiveDocument

# --- Date heading -----------------------------------------------------
$d.Content.Find.Execute("2024-12-30 Monday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2024-12-31 Tuesday", 2)

# --- Table of division problems ----------------------------------------
# The table has 5 "content" rows (1, 5, 9, 13, 17) with 5 cells each;
# several source expressions repeat (e.g. "75÷4=" appears twice, and
# "51÷7=" is both a source and a target elsewhere), so cells are
# addressed positionally rather than via a global text search.
$t = $d.Tables.Item(1)

$replacements = @(
    @(1, 1, "79÷5="),
    @(1, 2, "50÷4="),
    @(1, 3, "54÷3="),
    @(1, 4, "40÷6="),
    @(1, 5, "54÷7="),

    @(5, 1, "47÷7="),
    @(5, 2, "32÷3="),
    @(5, 3, "38÷4="),
    @(5, 4, "27÷6="),
    @(5, 5, "75÷2="),

    @(9, 1, "59÷8="),
    @(9, 2, "26÷7="),
    @(9, 3, "47÷5="),
    @(9, 4, "66÷5="),
    @(9, 5, "74÷9="),

    @(13, 1, "91÷9="),
    @(13, 2, "46÷8="),
    @(13, 3, "89÷7="),
    @(13, 4, "28÷3="),
    @(13, 5, "44÷5="),

    @(17, 1, "41÷8="),
    @(17, 2, "38÷6="),
    @(17, 3, "58÷6="),
    @(17, 4, "51÷7="),
    @(17, 5, "78÷6=")
)

foreach ($entry in $replacements) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $text
}
